$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-set numeric values on the columns that must remain stored as numbers
# (must happen before the text-format (numFmtId 49) styling is copied in,
# otherwise the runtime infers a text type for them).
$ws.Range("F36").Value = 2
$ws.Range("F38").Value = 30
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 32

# Clone the formatting (fills/borders/number formats) of row 35 onto the
# five new rows 36-40.
$ws.Range("A35:F35").Copy()
$ws.Range("A36:F40").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Give the new rows their own highlight color (adds a new fill + cellXfs
# trio, mirroring the existing "dark" row style but in red).
$ws.Range("A36:F40").Interior.Color = 7039985

# Row 36
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "CII1000"
$ws.Range("C36").Value = "CONTABILIDAD Y COSTOS"
$ws.Range("D36").Value = "47"
$ws.Range("E36").Value = 7

# Row 37
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "CIT2005"
$ws.Range("C37").Value = "INGENIERÍA DE SOFTWARE"
$ws.Range("D37").Value = "40, 42"
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = "24, 25"

# Row 38
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = "CIT2102"
$ws.Range("C38").Value = "COMUNICACIONES DIGITALES"
$ws.Range("D38").Value = "42"
$ws.Range("E38").Value = 7

# Row 39
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = "FIC1003"
$ws.Range("C39").Value = "DERECHO EN INGENIERÍA"
$ws.Range("D39").Value = "54"
$ws.Range("E39").Value = 7

# Row 40
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = "CIT2104"
$ws.Range("C40").Value = "ARQUITECTURA DE COMPUTADORES"
$ws.Range("D40").Value = "54"
$ws.Range("E40").Value = 7

# Window / selection bookkeeping to mirror the final author state
[void]$ws.Range("G23").Select()
$excel.ActiveWindow.Zoom = 100
